$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on every cell we touch so the values round-trip
# as literal strings (matching the workbook's original inline-string cells)
# instead of being auto-coerced into numbers/percentages by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "256.22"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.32%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.12"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.88%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.626"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-11.10%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05896"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.44%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.642"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.83%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8684"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.24%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9474"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.02%"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "One"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.01041"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1,617.33%"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1402"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.79%"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03748"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "8.60%"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07090"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.09%"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03209"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.22%"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09256"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.47%"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitForexToken"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001548"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.27%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006043"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.83%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.514"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.47%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.194"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.55%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.202"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.05%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3074"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.24%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1282"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.02%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.848"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.92%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04234"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.23%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.65%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004274"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.42%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001201"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.01%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001508"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2.85%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03815"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.08%"

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "KickToken"

$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006224"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "10.70%"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "BKEXToken"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1099"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.26%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002256"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.80%"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.65%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005503"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.07%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.03%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06024"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-33.09%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002280"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "7.05%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
